$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the driver vintage date from E12 (previously "2022-08-29")
$ws.Range("E12").ClearContents()

# Update the total samples count for B14 (weekly driver report refresh)
$ws.Range("B14").Value = 265400
